$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G40").Value = "Fallo"
$ws.Range("H40").Value = -1

$ws.Range("G41").Value = "Fallo"
$ws.Range("H41").Value = -1

$ws.Range("G44").Value = "Fallo"
$ws.Range("H44").Value = -1

$ws.Range("G46").Value = "Acierto"
$ws.Range("H46").Value = 1

$ws.Range("G47").Value = "Acierto"
$ws.Range("H47").Value = 0.57
